$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 32 and 33 swap their species-record data. Only touch the specific
# cells that actually change between the two rows (leave shared /
# unrelated cells such as dates, location, observer, etc. untouched).

$cols = @("A", "B", "E", "F", "G", "H", "M", "Q", "R")

foreach ($col in $cols) {
    $cell32 = $ws.Range("$col" + "32")
    $cell33 = $ws.Range("$col" + "33")

    $v32 = $cell32.Value()
    $v33 = $cell33.Value()

    $cell32.Value = $v33
    $cell33.Value = $v32
}

# "Publik kommentar" moves from row 32 to row 33.
$ws.Range("AC33").Value = $ws.Range("AC32").Value()
$ws.Range("AC32").ClearContents()
